# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" quarterly sheet (cloned from the existing
# "2022-Q2" sheet's layout/formatting) right after the "总计" summary sheet,
# fills it with the new quarter's fund data, and prepends a corresponding
# row to the "总计" summary sheet (shifting the existing history rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q3, push the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert(-4121, 0)   # xlShiftDown, xlFormatFromLeftOrAbove

# Copy row formatting from the (now shifted) old row2 = new row3 so the
# freshly inserted blank row matches the rest of the table's styling.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.05

# ---------------------------------------------------------------------
# 2) Quarterly sheets: clone "2022-Q2" (current tab #2) into a new sheet
#    placed right before it, rename it to "2022-Q3", and fill in the new
#    quarter's numbers.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

$q3Sheet.Range("D2:G2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "0.66"
$q3Sheet.Range("E2").Value = "86.08"
$q3Sheet.Range("F2").Value = "7.94"
$q3Sheet.Range("G2").Value = "0.0524"
$q3Sheet.Range("H2").Value = 2

# ---------------------------------------------------------------------
# 3) Restore the active/selected tab to the last sheet ("2020-Q4"), which
#    was the selected tab before this edit (the sheet-copy above makes
#    itself active, so move selection back).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
